# Updated Week numbers on time sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the week's dates back by 7 days (row 5, columns B:H) ---
$ws.Range("B5").Value = 44452
$ws.Range("C5").Value = 44453
$ws.Range("D5").Value = 44454
$ws.Range("E5").Value = 44455
$ws.Range("F5").Value = 44456
$ws.Range("G5").Value = 44457
$ws.Range("H5").Value = 44458

# --- Row 6 (Lecture): daily-total formula collapses to a plain value ---
$ws.Range("I6").Value = 1

# --- Row 8 (Team Meeting): hours moved onto Wednesday, daily total collapses to a plain value ---
$ws.Range("D8").Value = 1
$ws.Range("I8").Value = 1

# --- Row 9 (Sponsor Meeting): daily-total formula collapses to a plain value ---
$ws.Range("I9").Value = 1

# --- Row 10 (Organizing misc.): hour moved from Monday to Sunday, daily total collapses to a plain value ---
$ws.Range("B10").ClearContents()
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1

# --- Update selection / active cell ---
$ws.Range("G16").Select() | Out-Null
